$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Pre-Optimization table: "Division U9 (Tier: 0)" row (row 9) got rescheduled
# ---------------------------------------------------------------------------
$ws.Range("C9").Value = 16.0
$ws.Range("D9").Value = 6.0
$ws.Range("F9").Value = 16.0
$ws.Range("G9").Value = 6.0

# ---------------------------------------------------------------------------
# Post-Optimization table updates
# ---------------------------------------------------------------------------
# "Division U7 (Tier: 1)" row (row 16)
$ws.Range("C16").Value = 38.0
$ws.Range("D16").Value = 6.0
$ws.Range("F16").Value = 35.0
$ws.Range("G16").Value = 9.0

# "Division U7 (Tier: 2)" row (row 17)
$ws.Range("C17").Value = 36.0
$ws.Range("D17").Value = 8.0

# "Division U8 (Tier: 0)" row (row 19)
$ws.Range("F19").Value = 16.0
$ws.Range("G19").Value = 6.0

# "Division U8 (Tier: 1)" row (row 20)
$ws.Range("F20").Value = 11.0
$ws.Range("G20").Value = 11.0

# "Division U8 (Tier: 2)" row (row 21)
$ws.Range("F21").Value = 36.0
$ws.Range("G21").Value = 8.0

# "Division U8 (Tier: 3)" row (row 22)
$ws.Range("C22").Value = 43.0
$ws.Range("D22").Value = 1.0
$ws.Range("F22").Value = 37.0
$ws.Range("G22").Value = 7.0

# "Division U9 (Tier: 0)" row (row 23)
$ws.Range("F23").Value = 16.0
$ws.Range("G23").Value = 6.0

# "Division U9 (Tier: 1)" row (row 24)
$ws.Range("F24").Value = 35.0
$ws.Range("G24").Value = 9.0

# "Division U9 (Tier: 2)" row (row 25)
$ws.Range("C25").Value = 37.0
$ws.Range("D25").Value = 7.0

# "Division U9 (Tier: 3)" row (row 26)
$ws.Range("C26").Value = 113.0
$ws.Range("D26").Value = 19.0
$ws.Range("F26").Value = 94.0
$ws.Range("G26").Value = 38.0

# ---------------------------------------------------------------------------
# Update the "Scheduling Success %" text cells (column H). These are stored
# as literal text (e.g. "72.73%"), not as percentage-formatted numbers, so a
# plain .Value assignment (which Excel auto-converts into a numeric percent)
# would change both the stored type and the cell style. Instead, write the
# text via a formula and then paste-special as values to collapse it back
# down to a literal shared string while leaving style/number-format intact.
# ---------------------------------------------------------------------------
function Set-TextValue($rangeAddress, $text) {
    $cell = $ws.Range($rangeAddress)
    $cell.Formula = "=""" + $text + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

Set-TextValue "H9" "72.73%"
Set-TextValue "H16" "86.36%"
Set-TextValue "H17" "81.82%"
Set-TextValue "H22" "97.73%"
Set-TextValue "H25" "84.09%"
Set-TextValue "H26" "85.61%"

$excel.CutCopyMode = 0
